# Hana_T303.xlsx - "Adding the method to fetch new customer with phone number"
#
# The customer address (custaddress1, column B row 2) is updated to a
# fuller address that now also implies the phone-number lookup flow.
# custzip/custphone (columns C/D) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the customer address value -------------------------------
# Setting .Value directly resets the cell's style (drops the inherited
# "quote prefix" formatting flag picked up from the row). Restore the
# original look by copying formats back from a neighboring cell that
# still carries the same style after the value write.
$ws.Range("B2").Value = "2715 35th Avenue Greeley, CO, USA"
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# --- Widen column B so the longer address is readable -----------------
$ws.Range("B1").EntireColumn.ColumnWidth = 30.7265625

# --- Move the active selection to A3 (matches the saved view state) ---
$ws.Range("A3").Select()
